# Ajustando para calcular a correlação
# Clears the old "Correlação" label + CORREL(...) result from row 14 so the
# sheet can be recalculated fresh: unmerge B14:C14, wipe the three cells,
# and drop the thin border that used to frame that summary row (keeping the
# same fonts/fills/number-formats the row already had).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The summary cells used to be merged (B14:C14) - split them back apart
# before touching their contents.
$ws.Range("B14:C14").UnMerge()

# Remove the "Correlação" text label and the =CORREL(...) formula/result.
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()

# Row 14 no longer has the box border around it.
$ws.Range("A14").Borders.LineStyle = -4142
$ws.Range("B14:C14").Borders.LineStyle = -4142

# B14/C14 were center-aligned horizontally as part of the merged summary
# cell; now that it's just a plain (empty) pair of cells, drop the explicit
# horizontal alignment and keep only the vertical centering.
$ws.Range("B14:C14").HorizontalAlignment = 1

# Leave the cursor where the author last left it when they saved.
$ws.Range("E13").Select() | Out-Null
